# Applies the "add few fields and error report" update:
#  - row 2 (GS WORK O-RING GASKET / 5116877-50227802194) is replaced by the
#    NYLON RUBBER BAND Pkt 500gm / 5116877-32863368197 catalog entry and
#    gains an "Already exits" status in column X
#  - row 3 (Workstore 100202050014 / 5116877-37196116890) is replaced by the
#    GOOD MAKE RUBBER BAND / 5116877-15821908934 catalog entry
#  - row 4 (Workstore 100202050015 / 5116877-54522434970) is replaced by the
#    R-73947 / 5116877-17154744803 catalog entry and gains an
#    "Already exits" status in column X
#  - a brand-new row 5 is appended for RUBBER BAND GOOD MAKE / 5116877-92981260387

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The engine only ever *appends* hyperlink entries (editing .Address/.SubAddress
# in place, or re-Add()-ing over an existing anchor, just tacks on a duplicate),
# so the reliable way to rewrite them is to drop every hyperlink up front and
# recreate the full set afterwards, in the right left-to-right / top-to-bottom
# order.
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# Row 2: 5116877-32863368197 - NYLON RUBBER BAND Pkt 500gm
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "5116877-32863368197"
$ws.Range("B2").Value = "NYLON RUBBER BAND Pkt 500gm"
$ws.Range("D2").Value = "Ekam"
$ws.Range("O2").Value = "312.0"
$ws.Range("P2").Value = "279.99"
$ws.Range("T2").Value = "51415.0"
$ws.Range("U2").Value = "6.0"
$ws.Range("X2").Value = "Already exits"

# ---------------------------------------------------------------------
# Row 3: 5116877-15821908934 - GOOD MAKE RUBBER BAND
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "5116877-15821908934"
$ws.Range("B3").Value = "GOOD MAKE RUBBER BAND"
$ws.Range("D3").Value = "Good Make"
$ws.Range("O3").Value = "480.0"
$ws.Range("P3").Value = "398.99"
$ws.Range("T3").Value = "51417.0"
$ws.Range("U3").Value = "4.0"
$ws.Range("X3").Value = "Published"

# ---------------------------------------------------------------------
# Row 4: 5116877-17154744803 - R-73947
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "5116877-17154744803"
$ws.Range("B4").Value = "R-73947"
$ws.Range("D4").Value = "Good Make"
$ws.Range("O4").Value = "430.0"
$ws.Range("P4").Value = "380.99"
$ws.Range("T4").Value = "51418.0"
$ws.Range("U4").Value = "4.0"
$ws.Range("X4").Value = "Already exits"

# ---------------------------------------------------------------------
# Row 5 (new): 5116877-92981260387 - RUBBER BAND GOOD MAKE
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "5116877-92981260387"
$ws.Range("B5").Value = "RUBBER BAND GOOD MAKE"
$ws.Range("C5").Value = "rubber bands"
$ws.Range("D5").Value = "Good Make"
$ws.Range("G5").Value = "ST-GI12/A472-20-25"
$ws.Range("H5").Value = "Genaric"
$ws.Range("I5").Value = "44166.0"
$ws.Range("J5").Value = "44166.0"
$ws.Range("K5").Value = "46022.0"
$ws.Range("L5").Value = "India"
$ws.Range("N5").Value = "4016.0"
$ws.Range("O5").Value = "430.0"
$ws.Range("P5").Value = "380.99"
$ws.Range("S5").Value = "all"
$ws.Range("T5").Value = "51420.0"
$ws.Range("U5").Value = "4.0"
$ws.Range("V5").Value = "3.0"
$ws.Range("X5").Value = "Already exits"

# ---------------------------------------------------------------------
# Recreate the hyperlinks, left-to-right / top-to-bottom, so rId1..rId12
# come out in the same order as the target workbook.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E2"), "http://mkp.gem.gov.in/rubber-bands/nylon-rubber-band-pkt-500gm/p-5116877-32863368197-cat.html", "")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://admin-mkp.gem.gov.in/", "!/catalog/new?bnid=home_offi_of45811733_fast_rubb&gem_catalog_id=5116877-32863368197")
$ws.Hyperlinks.Add($ws.Range("W2"), "https://admin-mkp.gem.gov.in/", "!/catalog/new?bnid=home_offi_of45811733_fast_rubb&gem_catalog_id=5116877-32863368197")

$ws.Hyperlinks.Add($ws.Range("E3"), "http://mkp.gem.gov.in/rubber-bands/rubber-band-big-size/p-5116877-15821908934-cat.html", "")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://admin-mkp.gem.gov.in/", "!/catalog/new?bnid=home_offi_of45811733_fast_rubb&gem_catalog_id=5116877-15821908934")
$ws.Hyperlinks.Add($ws.Range("W3"), "https://admin-mkp.gem.gov.in/", "!/catalog/new?id=2382224-19041754072-cat&bnid=home_offi_of45811733_fast_rubb&gem_catalog_id=5116877-15821908934")

$ws.Hyperlinks.Add($ws.Range("E4"), "http://mkp.gem.gov.in/office-equipment-accessories-supplies/rubber-bands/p-5116877-17154744803-cat.html", "")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://admin-mkp.gem.gov.in/", "!/catalog/new?bnid=home_offi_of45811733_fast_rubb&gem_catalog_id=5116877-17154744803")
$ws.Hyperlinks.Add($ws.Range("W4"), "https://admin-mkp.gem.gov.in/", "!/catalog/new?bnid=home_offi_of45811733_fast_rubb&gem_catalog_id=5116877-17154744803")

$ws.Hyperlinks.Add($ws.Range("E5"), "http://mkp.gem.gov.in/rubber-bands/rubber-band-good-make/p-5116877-92981260387-cat.html", "")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://admin-mkp.gem.gov.in/", "!/catalog/new?bnid=home_offi_of45811733_fast_rubb&gem_catalog_id=5116877-92981260387")
$ws.Hyperlinks.Add($ws.Range("W5"), "https://admin-mkp.gem.gov.in/", "!/catalog/new?id=2382224-42904541515-cat&bnid=home_offi_of45811733_fast_rubb&gem_catalog_id=5116877-92981260387")

# Re-adding hyperlinks stacks a fresh format on top of whatever the cell had;
# pin every hyperlink cell back to the workbook's single "Hyperlink" cell
# style so we don't pick up extra, unused style variants.
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("W2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("W3").Style = "Hyperlink"
$ws.Range("E4").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("W4").Style = "Hyperlink"
$ws.Range("E5").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("W5").Style = "Hyperlink"

Write-Host ("Hyperlinks: " + $ws.Hyperlinks.Count() + ", dimension: " + $ws.UsedRange.Address())
